# Update cryptocurrency price/volume data (and two pairs of swapped rank rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='62.156.37'; E='  +3.12%  ' },
    @{ Row=3; D='2.680.52'; E='  +3.00%  ' },
    @{ Row=4; D='0.998'; E='  -0.20%  ' },
    @{ Row=5; D='585.37'; E='  +1.68%  ' },
    @{ Row=6; D='146.52'; E='  +2.48%  ' },
    @{ Row=7; D='0.996'; E='  -0.28%  ' },
    @{ Row=8; D='0.603'; E='  +0.77%  ' },
    @{ Row=9; D='6.62'; E='  +0.88%  ' },
    @{ Row=10; D='0.112'; E='  +5.81%  ' },
    @{ Row=11; D='0.384'; E='  +3.92%  ' },
    @{ Row=12; E='  +1.20%  ' },
    @{ Row=13; D='3.140.43'; E='  +2.67%  ' },
    @{ Row=14; D='26.25'; E='  +8.02%  ' },
    @{ Row=15; D='62.031.05'; E='  +2.93%  ' },
    @{ Row=16; D='0.0000149'; E='  +5.46%  ' },
    @{ Row=17; D='2.677.07'; E='  +2.80%  ' },
    @{ Row=18; D='11.79'; E='  +3.71%  ' },
    @{ Row=19; D='4.84'; E='  +4.64%  ' },
    @{ Row=20; D='360.20'; E='  +4.00%  ' },
    @{ Row=21; D='6.97'; E='  +1.18%  ' },
    @{ Row=22; D='1.00'; E='  +0.09%  ' },
    @{ Row=23; D='0.528'; E='  -0.19%  ' },
    @{ Row=24; D='65.07'; E='  +3.10%  ' },
    @{ Row=25; E='  +3.61%  ' },
    @{ Row=26; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='8.50'; E='  +6.07%  ' },
    @{ Row=27; B='Binance-PegBSC-USD'; C='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D='0.997'; E='  -0.22%  ' },
    @{ Row=28; E='  +8.40%  ' },
    @{ Row=29; D='0.0₃0839'; E='  +5.31%  ' },
    @{ Row=30; D='7.04'; E='  +10.43%  ' },
    @{ Row=31; D='170.01'; E='  +2.24%  ' },
    @{ Row=32; D='0.997'; E='  -0.08%  ' },
    @{ Row=33; D='20.40'; E='  +5.06%  ' },
    @{ Row=34; D='1.15'; E='  +16.46%  ' },
    @{ Row=35; D='4.72'; E='  +10.42%  ' },
    @{ Row=36; D='1.40'; E='  +7.67%  ' },
    @{ Row=37; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.77'; E='  +8.87%  ' },
    @{ Row=38; B='SuiNetwork'; C='https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; D='1.01'; E='  +20.09%  ' },
    @{ Row=39; D='348.85'; E='  +11.41%  ' },
    @{ Row=40; D='4.20'; E='  +7.90%  ' },
    @{ Row=41; D='38.70'; E='  +1.60%  ' },
    @{ Row=42; D='5.53'; E='  +10.75%  ' },
    @{ Row=43; D='21.21'; E='  +6.91%  ' },
    @{ Row=44; D='0.0586'; E='  +6.31%  ' },
    @{ Row=45; D='21.48'; E='  +7.51%  ' },
    @{ Row=46; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='136.72'; E='  +1.06%  ' },
    @{ Row=47; D='0.634'; E='  +4.79%  ' },
    @{ Row=48; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0256'; E='  +6.15%  ' },
    @{ Row=49; D='0.101'; E='  +1.31%  ' },
    @{ Row=50; D='0.993'; E='  -0.58%  ' },
    @{ Row=51; D='2.127.92'; E='  +5.69%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force text format so numeric-looking price strings (e.g. "1.00", "61.720.68")
        # keep their exact original formatting instead of being parsed as numbers
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
}
